# ValueSet-us-core-clinical-note-type.xlsx
# Updates the "Metadata" sheet: Publisher/Contact are changed to the new
# HL7 International / Cross-Group Projects contact, the previous
# "Health eData Inc" contact is kept as an additional Contact row, and
# Jurisdiction is now populated with "United States of America". All rows
# below the insertion point shift down by one.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Capture the current (pre-edit) text of the rows that are about to shift
# down, before any writes happen.
$oldA11 = $ws1.Range("A11").Text   # "Jurisdiction"
$oldA12 = $ws1.Range("A12").Text   # "Description"
$oldB12 = $ws1.Range("B12").Text
$oldA13 = $ws1.Range("A13").Text   # "Purpose"
$oldB13 = $ws1.Range("B13").Text
$oldA14 = $ws1.Range("A14").Text   # "Copyright"
$oldB14 = $ws1.Range("B14").Text
$oldA15 = $ws1.Range("A15").Text   # "Immutable"
$oldB15 = $ws1.Range("B15").Text

# Row 9: Publisher value
$ws1.Range("B9").Value = "HL7 International / Cross-Group Projects"

# Row 10: Contact value
$ws1.Range("B10").Value = "HL7 International / Cross-Group Projects (http://www.hl7.org/Special/committees/cgp, cgp@lists.HL7.org)"

# Row 11 (new): second Contact row, preserving the former publisher's contact
$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "Health eData Inc (mailto:ehaas@healthedatainc.com)"

# Row 12: Jurisdiction, now with a value
$ws1.Range("A12").Value = $oldA11
$ws1.Range("B12").Value = "United States of America"

# Row 13: Description (shifted down from row 12)
$ws1.Range("A13").Value = $oldA12
$ws1.Range("B13").Value = $oldB12

# Row 14: Purpose (shifted down from row 13)
$ws1.Range("A14").Value = $oldA13
$ws1.Range("B14").Value = $oldB13

# Row 15: Copyright (shifted down from row 14)
$ws1.Range("A15").Value = $oldA14
$ws1.Range("B15").Value = $oldB14

# Row 16 (new last row): Immutable (shifted down from row 15)
$ws1.Range("A16").Value = $oldA15
$ws1.Range("B16").Value = $oldB15
